# Apply the changes described by the commit diff:
#  - credentials sheet: update row-2 selections (C2, D2, F2), add a new
#    value in A4, and move the active selection to E4.
#  - halka sheet: move active selection to J22.
#  - villages sheet: move active selection to D9.
#  - delete sheet: move active selection to E7.

$wb = $excel.ActiveWorkbook

# --- credentials sheet ---------------------------------------------------
$credentials = $wb.Worksheets.Item("credentials")

$credentials.Range("C2").Value = "116374/हरनगला"
$credentials.Range("D2").Value = 600
$credentials.Range("F2").Value = "jayad"
$credentials.Range("A4").Value = "https://www.youtube.com/watch?v=gmfMUZyoGXg"

$credentials.Activate()
$credentials.Range("E4").Select()

# --- halka sheet -----------------------------------------------------------
$halka = $wb.Worksheets.Item("halka")
$halka.Activate()
$halka.Range("J22").Select()

# --- villages sheet ----------------------------------------------------
$villages = $wb.Worksheets.Item("villages")
$villages.Activate()
$villages.Range("D9").Select()

# --- delete sheet --------------------------------------------------------
$delete = $wb.Worksheets.Item("delete")
$delete.Activate()
$delete.Range("E7").Select()

# restore original active sheet (credentials was tabSelected in the source)
$credentials.Activate()
